# Chapter 6 example 6: add the "by region / by division" summary table
# (J1:O6) next to the existing 6-month sales detail table, mirroring a
# PivotTable of Sum(销售金额) grouped by 销售地区 (rows) x 销售分部 (columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ------------------------------------------------
$ws.Range("J1").Value = "销售地区"
$ws.Range("K1").Value = "销售一分部"
$ws.Range("L1").Value = "销售三分部"
$ws.Range("M1").Value = "销售二分部"
$ws.Range("N1").Value = "销售四分部"
$ws.Range("O1").Value = "总计"

# ---- Row labels (column J) ---------------------------------------------
$ws.Range("J2").Value = "华东"
$ws.Range("J3").Value = "华中"
$ws.Range("J4").Value = "华北"
$ws.Range("J5").Value = "华南"
$ws.Range("J6").Value = "总计"

# ---- Values --------------------------------------------------------------
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 2059200
$ws.Range("M2").Value = 4183800
$ws.Range("N2").Value = 3513200
$ws.Range("O2").Value = 9756200

$ws.Range("K3").Value = 3826000
$ws.Range("L3").Value = 1806200
$ws.Range("M3").Value = 6324000
$ws.Range("N3").Value = 1284000
$ws.Range("O3").Value = 13240200

$ws.Range("K4").Value = 3676400
$ws.Range("L4").Value = 1694000
$ws.Range("M4").Value = 1245200
$ws.Range("N4").Value = 2552800
$ws.Range("O4").Value = 9168400

$ws.Range("K5").Value = 3025200
$ws.Range("L5").Value = 1634600
$ws.Range("M5").Value = 588000
$ws.Range("N5").Value = 3369400
$ws.Range("O5").Value = 8617200

$ws.Range("K6").Value = 10527600
$ws.Range("L6").Value = 7194000
$ws.Range("M6").Value = 12341000
$ws.Range("N6").Value = 10719400
$ws.Range("O6").Value = 40782000

# ---- Number formatting: currency with 2 decimals for the money cells ----
$moneyFormat = """¥""#,##0.00;""¥""\-#,##0.00"
$ws.Range("L2:O6").NumberFormat = $moneyFormat

# ---- Center alignment to match the rest of the summary block -----------
$ws.Range("J1:O6").HorizontalAlignment = -4108
$ws.Range("J1:O6").VerticalAlignment = -4108

# ---- Autofit columns, matching Excel's behaviour after entering data ---
$ws.Range("A1:O42").Columns.AutoFit()
